$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "null"
$ws.Range("A3").Value = "PEP_ID-2010444"
